$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, shifting existing rows 9-15 down to 10-16.
$ws.Rows.Item(9).Insert()

# Populate the new row 9 with the latest weekly price record.
$ws.Cells.Item(9, 1).Value = 1
$ws.Cells.Item(9, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(9, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(9, 4).Value = 44607
$ws.Cells.Item(9, 5).Value = 15
$ws.Cells.Item(9, 6).Value = 100112026
$ws.Cells.Item(9, 7).Value = "Haba"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 900
$ws.Cells.Item(9, 11).Value = 1300
$ws.Cells.Item(9, 12).Value = 1400
$ws.Cells.Item(9, 13).Value = 1350
$ws.Cells.Item(9, 14).Value = "$/kilo"
$ws.Cells.Item(9, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(9, 16).Value = 1350
$ws.Cells.Item(9, 17).Value = 1
$ws.Cells.Item(9, 18).Value = "Hortaliza"
